$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "26.250.52"
$ws.Range("E2").Value2 = "  -0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.657.77"
$ws.Range("E3").Value2 = "  -0.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.004"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "219.61"
$ws.Range("E5").Value2 = "  -0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.5238"
$ws.Range("E6").Value2 = "  -2.07%  "
$ws.Range("E7").Value2 = "  -0.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.2674"
$ws.Range("E8").Value2 = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.06360"
$ws.Range("E9").Value2 = "  -0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "20.72"
$ws.Range("E10").Value2 = "  -1.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.07771"
$ws.Range("E11").Value2 = "  -0.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "4.566"
$ws.Range("E12").Value2 = "  +0.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "1.660.83"
$ws.Range("E13").Value2 = "  -0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "1.885.49"
$ws.Range("E14").Value2 = "  -0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.5694"
$ws.Range("E15").Value2 = "  +0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.0₅8169"
$ws.Range("E16").Value2 = "  -0.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "65.57"
$ws.Range("E17").Value2 = "  -1.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "26.240.31"
$ws.Range("E18").Value2 = "  -0.85%  "
$ws.Range("E19").Value2 = "  -0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "4.722"
$ws.Range("E20").Value2 = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "192.23"
$ws.Range("E21").Value2 = "  -3.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "10.38"
$ws.Range("E22").Value2 = "  +0.14%  "
$ws.Range("E23").Value2 = "  -0.59%  "
$ws.Range("E24").Value2 = "  -0.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "143.62"
$ws.Range("E25").Value2 = "  -2.14%  "
$ws.Range("E26").Value2 = "  -2.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "7.281"
$ws.Range("E27").Value2 = "  +0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "16.01"
$ws.Range("E28").Value2 = "  -1.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "1.492"
$ws.Range("E29").Value2 = "  -0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0.05644"
$ws.Range("E30").Value2 = "  -4.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "1.277"
$ws.Range("E31").Value2 = "  -0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "3.504"
$ws.Range("E32").Value2 = "  -2.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "3.380"
$ws.Range("E33").Value2 = "  +1.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "1.587"
$ws.Range("E34").Value2 = "  -2.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "2.805"
$ws.Range("E35").Value2 = "  -1.71%  "
$ws.Range("E36").Value2 = "  -2.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "2.407"
$ws.Range("E37").Value2 = "  -1.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.5778"
$ws.Range("E38").Value2 = "  -1.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.01601"
$ws.Range("E39").Value2 = "  -1.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "5.909"
$ws.Range("E40").Value2 = "  -0.31%  "
$ws.Range("E41").Value2 = "  -0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.8486"
$ws.Range("E42").Value2 = "  -2.15%  "
$ws.Range("E43").Value2 = "  -0.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "1.032.08"
$ws.Range("E44").Value2 = "  -4.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "102.30"
$ws.Range("E45").Value2 = "  -1.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "1.795.85"
$ws.Range("E46").Value2 = "  -1.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "58.65"
$ws.Range("E47").Value2 = "  +0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.0₈105"
$ws.Range("E48").Value2 = "  -1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "1.004"
$ws.Range("E49").Value2 = "  -1.01%  "
$ws.Range("E50").Value2 = "  +2.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "8.029"
$ws.Range("E51").Value2 = "  -0.36%  "
